# [Resource Manager] Fix header of template file
#
# The "Host Information" header row documents the allowed values for the
# "Disks" column. Update the inline note from
#   "[status: 0 - available, 1 - inused]"
# to
#   "[status(optional): 0 - available, 1 - reserved]"
# while preserving the existing rich-text formatting (the note is bold,
# unlike the "Disks" line above it), then leave the sheet's selection on
# O2 (matching the refreshed template's last-saved cursor position)
# instead of Q1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Host Information")

# --- Update the Q1 header note text, run-by-run, so the "Disks" run keeps
#     its own (non-bold) formatting and only the bracketed note run's text
#     changes. ---
$cell = $ws.Range("Q1")

$oldNote = "[status: 0 - available, 1 - inused]"
$newNote = "[status(optional): 0 - available, 1 - reserved]"

$fullText = $cell.Characters().Text
$noteStart = $fullText.IndexOf($oldNote)
if ($noteStart -ge 0) {
    # Characters() is 1-based.
    $noteChars = $cell.Characters($noteStart + 1, $oldNote.Length)
    $noteChars.Text = $newNote

    # Re-apply the bold note formatting to the replaced run so it keeps
    # looking like the rest of the header note (font name/size/weight).
    $newNoteChars = $cell.Characters($noteStart + 1, $newNote.Length)
    $newFont = $newNoteChars.Font
    $newFont.Bold = $true
    $newFont.Size = 10
    $newFont.Name = "宋体"
}

# --- Move the sheet's saved selection from Q1 to O2. ---
$ws.Range("O2").Select()
